# Auto-generated Excel COM-interop script implementing the commit diff.
# Adds a new randomsearch result row and refreshes downstream metrics/selector reprs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long-form cell text templates ---
$b2Text = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5c8cc1ba60>),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])
"@
$c2Text = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f5e8c832610>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}
"@
$b3Text = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5e8c832580>),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])
"@
$c3Text = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f5e8c84d220>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}
"@
$bPlainText = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])
"@
$cPlainText = @"
{'selector': RandomUnderSampler(random_state=42), 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}
"@
$b6Text = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5e8cb2fc40>),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])
"@
$c6Text = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f5e8cb2fb20>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}
"@

# --- Row structure: insert a fresh result row at 3, drop the row that is
#     superseded (the old row 5 once its content has shifted to row 6) ---
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(6).Delete()

# Re-apply the bordered/centered header-row style to the new A3 cell
# (row insert clones formatting from the row above, but not the exact style id).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2: selector swapped for the NamedFeatureSelector instance, metrics refreshed ---
$ws.Range("B2").Value = $b2Text
$ws.Range("C2").Value = $c2Text
$ws.Rows.Item(2).AutoFit()
$ws.Range("D2").Value = 0.5405139563626574
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 0.6863482744004518
$ws.Range("H2").Value = 0.4085081585081585
$ws.Range("I2").Value = "[1 1 1 0 1 0 1 0 1 0 1 1 1 0 0 0 0 1 0 0 0 0 1 0]"
$ws.Range("J2").Value = "[0 1 0 1 0 0 1 1 1 1 0 1 1 0 1 1 0 1 1 1 1 1 0 0]"

# --- Row 3 (new): freshly inserted search result ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = $b3Text
$ws.Range("C3").Value = $c3Text
$ws.Rows.Item(3).AutoFit()
$ws.Range("D3").Value = 0.5332457799226544
$ws.Range("E3").Value = "Random"
$ws.Range("F3").Value = 69
$ws.Range("G3").Value = 0.6857825567502988
$ws.Range("H3").Value = 0.2879710144927536
$ws.Range("I3").Value = "[0 1 1 0 1 0 0 0 1 1 1 0 0 0 1 0 1 0 1 1 0 0 1 0]"
$ws.Range("J3").Value = "[1 1 0 0 0 1 1 1 0 1 1 1 1 1 0 1 0 0 0 0 0 1 1 1]"

# --- Row 4 (previously row 3): selector untouched, metrics refreshed ---
$ws.Range("D4").Value = 0.5902163898155958
$ws.Range("F4").Value = 23
$ws.Range("G4").Value = 0.7136080670885425
$ws.Range("H4").Value = 0.5648148148148149
$ws.Range("I4").Value = "[0 1 1 0 0 1 0 0 0 1 0 1 0 1 0 1 0 0 1 0 0 1 1 1]"
$ws.Range("J4").Value = "[0 0 0 1 1 1 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 1 0]"

# --- Row 5 (previously row 4): selector untouched, metrics refreshed ---
$ws.Range("D5").Value = 0.4858341144753852
$ws.Range("F5").Value = 42
$ws.Range("G5").Value = 0.7307474641939838
$ws.Range("H5").Value = 0.3760869565217391
$ws.Range("I5").Value = "[1 0 1 0 0 0 1 1 1 1 1 1 1 0 0 0 0 0 0 0 1 1 0 0]"
$ws.Range("J5").Value = "[0 0 0 0 0 1 0 0 0 0 1 1 1 0 0 1 1 1 1 1 1 0 1 1]"

# --- Row 6 (previously row 6): selector swapped for NamedFeatureSelector, metrics refreshed ---
$ws.Range("B6").Value = $b6Text
$ws.Range("C6").Value = $c6Text
$ws.Rows.Item(6).AutoFit()
$ws.Range("D6").Value = 0.564277385131134
$ws.Range("F6").Value = 89
$ws.Range("G6").Value = 0.6424422197503399
$ws.Range("H6").Value = 0.5833333333333334
$ws.Range("I6").Value = "[1 0 1 0 1 1 0 0 0 1 0 1 0 1 1 1 0 1 0 1 0 0 0 0]"
$ws.Range("J6").Value = "[0 1 1 0 1 0 0 0 1 0 0 1 0 0 1 1 0 1 1 1 1 0 1 1]"

